# Update market-board derived leve profit figures per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 869.939
$ws.Range("J129").Value = 959.5571
$ws.Range("L129").Value = 2878.6713
$ws.Range("N129").Value = -12878.6713
$ws.Range("H132").Value = 237298.1
$ws.Range("I132").Value = 4950.1816
$ws.Range("J132").Value = 1004046.2
$ws.Range("K132").Value = 14850.5448
$ws.Range("L132").Value = 3012138.6
$ws.Range("M132").Value = -12320.5448
$ws.Range("N132").Value = -3017198.6
$ws.Range("H135").Value = 495.6
$ws.Range("I135").Value = 388.14285
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 3493.28565
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -958.2856500000003
$ws.Range("N135").Value = -23070
$ws.Range("H137").Value = 2536.1396
$ws.Range("I137").Value = 1426.5217
$ws.Range("J137").Value = 3812.2
$ws.Range("K137").Value = 4279.5651
$ws.Range("L137").Value = 11436.6
$ws.Range("M137").Value = -1729.5651
$ws.Range("N137").Value = -16536.6
$ws.Range("H138").Value = 3217.82
$ws.Range("I138").Value = 1353
$ws.Range("J138").Value = 3379.9783
$ws.Range("K138").Value = 4059
$ws.Range("L138").Value = 10139.9349
$ws.Range("M138").Value = 1081
$ws.Range("N138").Value = -20419.9349
$ws.Range("H141").Value = 16957.215
$ws.Range("I141").Value = 20407.363
$ws.Range("J141").Value = 4306.6665
$ws.Range("K141").Value = 61222.08900000001
$ws.Range("L141").Value = 12919.9995
$ws.Range("M141").Value = -56042.08900000001
$ws.Range("N141").Value = -23279.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3797.4487
$ws.Range("I32").Value = 2964.6892
$ws.Range("K32").Value = 2964.6892
$ws.Range("M32").Value = -2677.6892
$ws.Range("H61").Value = 1380.7142
$ws.Range("I61").Value = 1046.6666
$ws.Range("J61").Value = 2215.8333
$ws.Range("K61").Value = 1046.6666
$ws.Range("L61").Value = 2215.8333
$ws.Range("M61").Value = -834.6666
$ws.Range("N61").Value = -2639.8333
$ws.Range("H74").Value = 2716.1555
$ws.Range("I74").Value = 2702.973
$ws.Range("K74").Value = 2702.973
$ws.Range("M74").Value = -1828.973
$ws.Range("H77").Value = 2716.1555
$ws.Range("I77").Value = 2702.973
$ws.Range("K77").Value = 13514.865
$ws.Range("M77").Value = -9146.865
$ws.Range("H80").Value = 32112.3
$ws.Range("J80").Value = 33569.223
$ws.Range("L80").Value = 33569.223
$ws.Range("N80").Value = -35565.223
$ws.Range("H83").Value = 32112.3
$ws.Range("J83").Value = 33569.223
$ws.Range("L83").Value = 100707.669
$ws.Range("N83").Value = -110691.669
$ws.Range("H132").Value = 2472.7632
$ws.Range("I132").Value = 1698.1154
$ws.Range("J132").Value = 4151.1665
$ws.Range("K132").Value = 5094.3462
$ws.Range("L132").Value = 12453.4995
$ws.Range("M132").Value = -2564.3462
$ws.Range("N132").Value = -17513.4995
$ws.Range("H136").Value = 1380.7142
$ws.Range("I136").Value = 1046.6666
$ws.Range("J136").Value = 2215.8333
$ws.Range("K136").Value = 3139.9998
$ws.Range("L136").Value = 6647.499899999999
$ws.Range("M136").Value = -589.9998000000001
$ws.Range("N136").Value = -11747.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 623.5
$ws.Range("I36").Value = 623.5
$ws.Range("K36").Value = 623.5
$ws.Range("M36").Value = -89.5
$ws.Range("H94").Value = 554.3913
$ws.Range("I94").Value = 657.2941
$ws.Range("J94").Value = 262.83334
$ws.Range("K94").Value = 657.2941
$ws.Range("L94").Value = 262.83334
$ws.Range("M94").Value = -206.2941
$ws.Range("N94").Value = -1164.83334
$ws.Range("H132").Value = 52988.5
$ws.Range("J132").Value = 52988.5
$ws.Range("L132").Value = 52988.5
$ws.Range("N132").Value = -63108.5
$ws.Range("H134").Value = 2421.6511
$ws.Range("I134").Value = 1233.6285
$ws.Range("J134").Value = 7619.25
$ws.Range("K134").Value = 3700.8855
$ws.Range("L134").Value = 22857.75
$ws.Range("M134").Value = -1165.8855
$ws.Range("N134").Value = -27927.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 30000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 30000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 30000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -30226
$ws.Range("H31").Value = 10871955
$ws.Range("I31").Value = 1212.6666
$ws.Range("J31").Value = 38466916
$ws.Range("K31").Value = 1212.6666
$ws.Range("L31").Value = 38466916
$ws.Range("M31").Value = -917.6666
$ws.Range("N31").Value = -38467506
$ws.Range("H34").Value = 10871955
$ws.Range("I34").Value = 1212.6666
$ws.Range("J34").Value = 38466916
$ws.Range("K34").Value = 1212.6666
$ws.Range("L34").Value = 38466916
$ws.Range("M34").Value = -1010.6666
$ws.Range("N34").Value = -38467320
$ws.Range("H58").Value = 1635.6824
$ws.Range("I58").Value = 1455.7028
$ws.Range("J58").Value = 2846.4546
$ws.Range("K58").Value = 1455.7028
$ws.Range("L58").Value = 2846.4546
$ws.Range("M58").Value = -1252.7028
$ws.Range("N58").Value = -3252.4546
$ws.Range("H134").Value = 5594.3447
$ws.Range("I134").Value = 7251.375
$ws.Range("K134").Value = 21754.125
$ws.Range("M134").Value = -19219.125
$ws.Range("H136").Value = 1635.6824
$ws.Range("I136").Value = 1455.7028
$ws.Range("J136").Value = 2846.4546
$ws.Range("K136").Value = 4367.1084
$ws.Range("L136").Value = 8539.363799999999
$ws.Range("M136").Value = -1817.1084
$ws.Range("N136").Value = -13639.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 925.75
$ws.Range("I57").Value = 452.5
$ws.Range("J57").Value = 1399
$ws.Range("K57").Value = 1357.5
$ws.Range("L57").Value = 4197
$ws.Range("M57").Value = -798.5
$ws.Range("N57").Value = -5315
$ws.Range("H131").Value = 8621553
$ws.Range("J131").Value = 895.86536
$ws.Range("L131").Value = 2687.59608
$ws.Range("N131").Value = -12767.59608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3129.64
$ws.Range("I132").Value = 1945.8125
$ws.Range("J132").Value = 5234.222
$ws.Range("K132").Value = 5837.4375
$ws.Range("L132").Value = 15702.666
$ws.Range("M132").Value = -3307.4375
$ws.Range("N132").Value = -20762.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22280.646
$ws.Range("I132").Value = 82966.664
$ws.Range("J132").Value = 9276.5
$ws.Range("K132").Value = 248899.992
$ws.Range("L132").Value = 27829.5
$ws.Range("M132").Value = -246369.992
$ws.Range("N132").Value = -32889.5
$ws.Range("H136").Value = 2336.15
$ws.Range("I136").Value = 1261.8485
$ws.Range("K136").Value = 3785.5455
$ws.Range("M136").Value = -1235.5455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 27781678
$ws.Range("I132").Value = 3150
$ws.Range("J132").Value = 33337384
$ws.Range("K132").Value = 9450
$ws.Range("L132").Value = 100012152
$ws.Range("M132").Value = -6920
$ws.Range("N132").Value = -100017212
$ws.Range("H136").Value = 1611.7969
$ws.Range("I136").Value = 528.24445
$ws.Range("J136").Value = 4178.1055
$ws.Range("K136").Value = 1584.73335
$ws.Range("L136").Value = 12534.3165
$ws.Range("M136").Value = 965.26665
$ws.Range("N136").Value = -17634.3165
